# Apply the updated metric values from the commit "BEST? done, now opti"
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Features")
$ws2 = $wb.Worksheets.Item("Global Metrics")

# --- Features sheet updates ---
$ws1.Range("B2:E2").NumberFormat = "@"
$ws1.Range("B2").Value = "0,694"
$ws1.Range("C2").Value = "0,833"
$ws1.Range("D2").Value = "0,758"
$ws1.Range("E2").Value = "0,990"

$ws1.Range("B3:E3").NumberFormat = "@"
$ws1.Range("B3").Value = "0,743"
$ws1.Range("C3").Value = "0,929"
$ws1.Range("D3").Value = "0,825"
$ws1.Range("E3").Value = "0,947"

$ws1.Range("B4:E4").NumberFormat = "@"
$ws1.Range("B4").Value = "0,639"
$ws1.Range("C4").Value = "0,742"
$ws1.Range("D4").Value = "0,687"
$ws1.Range("E4").Value = "0,971"

$ws1.Range("B5:E5").NumberFormat = "@"
$ws1.Range("B5").Value = "0,743"
$ws1.Range("C5").Value = "0,897"
$ws1.Range("D5").Value = "0,812"
$ws1.Range("E5").Value = "0,982"

$ws1.Range("B6:E6").NumberFormat = "@"
$ws1.Range("B6").Value = "0,281"
$ws1.Range("C6").Value = "0,290"
$ws1.Range("D6").Value = "0,286"
$ws1.Range("E6").Value = "0,991"

$ws1.Range("B7:E7").NumberFormat = "@"
$ws1.Range("B7").Value = "0,312"
$ws1.Range("C7").Value = "0,323"
$ws1.Range("D7").Value = "0,317"
$ws1.Range("E7").Value = "1,000"

$ws1.Range("B8:E8").NumberFormat = "@"
$ws1.Range("B8").Value = "0,185"
$ws1.Range("C8").Value = "0,385"
$ws1.Range("D8").Value = "0,250"
$ws1.Range("E8").Value = "0,956"

$ws1.Range("B9:E9").NumberFormat = "@"
$ws1.Range("B9").Value = "0,588"
$ws1.Range("C9").Value = "0,556"
$ws1.Range("D9").Value = "0,571"
$ws1.Range("E9").Value = "1,000"

$ws1.Range("B10:E10").NumberFormat = "@"
$ws1.Range("B10").Value = "0,375"
$ws1.Range("C10").Value = "0,562"
$ws1.Range("D10").Value = "0,450"
$ws1.Range("E10").Value = "0,913"

$ws1.Range("B11:E11").NumberFormat = "@"
$ws1.Range("B11").Value = "0,583"
$ws1.Range("C11").Value = "0,636"
$ws1.Range("D11").Value = "0,609"
$ws1.Range("E11").Value = "0,898"

$ws1.Range("B12:E12").NumberFormat = "@"
$ws1.Range("B12").Value = "0,583"
$ws1.Range("C12").Value = "0,636"
$ws1.Range("D12").Value = "0,609"
$ws1.Range("E12").Value = "0,918"

$ws1.Range("B13:E13").NumberFormat = "@"
$ws1.Range("B13").Value = "0,583"
$ws1.Range("C13").Value = "0,583"
$ws1.Range("D13").Value = "0,583"
$ws1.Range("E13").Value = "0,921"

$ws1.Range("B14:E14").NumberFormat = "@"
$ws1.Range("B14").Value = "0,500"
$ws1.Range("C14").Value = "0,545"
$ws1.Range("D14").Value = "0,522"
$ws1.Range("E14").Value = "0,858"

$ws1.Range("B15:E15").NumberFormat = "@"
$ws1.Range("B15").Value = "0,478"
$ws1.Range("C15").Value = "0,579"
$ws1.Range("D15").Value = "0,524"
$ws1.Range("E15").Value = "1,000"

$ws1.Range("B16:E16").NumberFormat = "@"
$ws1.Range("B16").Value = "0,375"
$ws1.Range("C16").Value = "0,643"
$ws1.Range("D16").Value = "0,474"
$ws1.Range("E16").Value = "0,971"

$ws1.Range("B17:E17").NumberFormat = "@"
$ws1.Range("B17").Value = "0,474"
$ws1.Range("C17").Value = "0,692"
$ws1.Range("D17").Value = "0,562"
$ws1.Range("E17").Value = "0,944"

$ws1.Range("B18:E18").NumberFormat = "@"
$ws1.Range("B18").Value = "0,400"
$ws1.Range("C18").Value = "0,571"
$ws1.Range("D18").Value = "0,471"
$ws1.Range("E18").Value = "0,949"

$ws1.Range("B19:E19").NumberFormat = "@"
$ws1.Range("B19").Value = "0,364"
$ws1.Range("C19").Value = "0,727"
$ws1.Range("D19").Value = "0,485"
$ws1.Range("E19").Value = "0,876"

$ws1.Range("B20:E20").NumberFormat = "@"
$ws1.Range("B20").Value = "0,364"
$ws1.Range("C20").Value = "0,727"
$ws1.Range("D20").Value = "0,485"
$ws1.Range("E20").Value = "0,876"

$ws1.Range("B21:E21").NumberFormat = "@"
$ws1.Range("B21").Value = "0,700"
$ws1.Range("C21").Value = "0,778"
$ws1.Range("D21").Value = "0,737"
$ws1.Range("E21").Value = "0,937"

$ws1.Range("B22:E22").NumberFormat = "@"
$ws1.Range("B22").Value = "0,500"
$ws1.Range("C22").Value = "0,300"
$ws1.Range("D22").Value = "0,375"
$ws1.Range("E22").Value = "1,000"

$ws1.Range("B23:E23").NumberFormat = "@"
$ws1.Range("B23").Value = "0,667"
$ws1.Range("C23").Value = "0,400"
$ws1.Range("D23").Value = "0,500"
$ws1.Range("E23").Value = "0,887"

$ws1.Range("B24:E24").NumberFormat = "@"
$ws1.Range("B24").Value = "0,143"
$ws1.Range("C24").Value = "0,091"
$ws1.Range("D24").Value = "0,111"
$ws1.Range("E24").Value = "0,707"

$ws1.Range("B25:E25").NumberFormat = "@"
$ws1.Range("B25").Value = "0,667"
$ws1.Range("C25").Value = "0,500"
$ws1.Range("D25").Value = "0,571"
$ws1.Range("E25").Value = "0,865"

$ws1.Range("B26:E26").NumberFormat = "@"
$ws1.Range("B26").Value = "0,154"
$ws1.Range("C26").Value = "0,500"
$ws1.Range("D26").Value = "0,235"
$ws1.Range("E26").Value = "1,000"

$ws1.Range("B28:E28").NumberFormat = "@"
$ws1.Range("B28").Value = "0,143"
$ws1.Range("C28").Value = "0,500"
$ws1.Range("D28").Value = "0,222"
$ws1.Range("E28").Value = "1,000"

$ws1.Range("B29:E29").NumberFormat = "@"
$ws1.Range("B29").Value = "0,105"
$ws1.Range("C29").Value = "0,250"
$ws1.Range("D29").Value = "0,148"
$ws1.Range("E29").Value = "1,000"

$ws1.Range("B30:E30").NumberFormat = "@"
$ws1.Range("B30").Value = "0,364"
$ws1.Range("C30").Value = "1,000"
$ws1.Range("D30").Value = "0,533"
$ws1.Range("E30").Value = "0,364"

$ws1.Range("B32:E32").NumberFormat = "@"
$ws1.Range("B32").Value = "0,133"
$ws1.Range("C32").Value = "0,333"
$ws1.Range("D32").Value = "0,190"
$ws1.Range("E32").Value = "1,000"

$ws1.Range("B33:E33").NumberFormat = "@"
$ws1.Range("B33").Value = "1,000"
$ws1.Range("C33").Value = "1,000"
$ws1.Range("D33").Value = "1,000"
$ws1.Range("E33").Value = "1,000"

$ws1.Range("E34").NumberFormat = "@"
$ws1.Range("E34").Value = "1,000"

$ws1.Range("B36:E36").NumberFormat = "@"
$ws1.Range("B36").Value = "0,714"
$ws1.Range("C36").Value = "0,455"
$ws1.Range("D36").Value = "0,556"
$ws1.Range("E36").Value = "1,000"

$ws1.Range("B37:E37").NumberFormat = "@"
$ws1.Range("B37").Value = "0,500"
$ws1.Range("C37").Value = "0,167"
$ws1.Range("D37").Value = "0,250"
$ws1.Range("E37").Value = "1,000"

$ws1.Range("E38").NumberFormat = "@"
$ws1.Range("E38").Value = "1,000"

$ws1.Range("E39").NumberFormat = "@"
$ws1.Range("E39").Value = "1,000"

$ws1.Range("B40:E40").NumberFormat = "@"
$ws1.Range("B40").Value = "0,667"
$ws1.Range("C40").Value = "0,500"
$ws1.Range("D40").Value = "0,571"
$ws1.Range("E40").Value = "1,000"

$ws1.Range("B41:E41").NumberFormat = "@"
$ws1.Range("B41").Value = "0,500"
$ws1.Range("C41").Value = "1,000"
$ws1.Range("D41").Value = "0,667"
$ws1.Range("E41").Value = "0,500"

$ws1.Range("E42").NumberFormat = "@"
$ws1.Range("E42").Value = "1,000"

$ws1.Range("B43:E43").NumberFormat = "@"
$ws1.Range("B43").Value = "0,500"
$ws1.Range("C43").Value = "0,400"
$ws1.Range("D43").Value = "0,444"
$ws1.Range("E43").Value = "1,000"

$ws1.Range("E45").NumberFormat = "@"
$ws1.Range("E45").Value = "1,000"

$ws1.Range("E46").NumberFormat = "@"
$ws1.Range("E46").Value = "1,000"

$ws1.Range("E47").NumberFormat = "@"
$ws1.Range("E47").Value = "1,000"

$ws1.Range("E48").NumberFormat = "@"
$ws1.Range("E48").Value = "1,000"

$ws1.Range("E53").NumberFormat = "@"
$ws1.Range("E53").Value = "0,475"

$ws1.Range("E58").NumberFormat = "@"
$ws1.Range("E58").Value = "1,000"

$ws1.Range("E61").NumberFormat = "@"
$ws1.Range("E61").Value = "0,379"

$ws1.Range("E62").NumberFormat = "@"
$ws1.Range("E62").Value = "1,000"

$ws1.Range("E63").NumberFormat = "@"
$ws1.Range("E63").Value = "1,000"

# --- Global Metrics sheet updates ---
$ws2.Range("B2:E2").NumberFormat = "@"
$ws2.Range("B2").Value = "0,290"
$ws2.Range("C2").Value = "0,735"
$ws2.Range("D2").Value = "0,467"
$ws2.Range("E2").Value = "0,896"

